# Refresh the cryptos price/volume table with latest scraped values
# (GitHub Actions cron update). A handful of price cells ("D" column)
# are plain decimal text like "563.65" / "1.00" that Excel would
# otherwise auto-coerce into numbers on assignment, so those are
# written with a leading apostrophe to force them to stay text, just
# like the original inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '62.278.06'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '2.430.64'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'" + '563.65'
$ws.Range('E5').Value = '  +0.32%  '
$ws.Range('D6').Value = "'" + '144.69'
$ws.Range('E6').Value = '  +0.71%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').Value = '2.429.87'
$ws.Range('E9').Value = '  +0.45%  '
$ws.Range('E10').Value = '  +0.28%  '
$ws.Range('D11').Value = "'" + '0.155'
$ws.Range('E11').Value = '  +0.36%  '
$ws.Range('D12').Value = "'" + '5.25'
$ws.Range('E12').Value = '  -1.74%  '
$ws.Range('D13').Value = "'" + '0.350'
$ws.Range('E13').Value = '  -0.97%  '
$ws.Range('D14').Value = "'" + '26.64'
$ws.Range('E14').Value = '  +2.99%  '
$ws.Range('E15').Value = '  -1.18%  '
$ws.Range('D16').Value = '2.843.68'
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('D17').Value = '62.221.39'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').Value = '2.427.45'
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('E19').Value = '  -0.39%  '
$ws.Range('D20').Value = "'" + '324.87'
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = "'" + '6.85'
$ws.Range('E21').Value = '  +0.41%  '
$ws.Range('B22').Value = 'Polkadot'
$ws.Range('C22').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D22').Value = "'" + '4.15'
$ws.Range('E22').Value = '  -1.36%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').Value = "'" + '67.39'
$ws.Range('E24').Value = '  +2.40%  '
$ws.Range('E25').Value = '  +0.92%  '
$ws.Range('D26').Value = "'" + '8.67'
$ws.Range('E26').Value = '  -3.89%  '
$ws.Range('D27').Value = "'" + '554.07'
$ws.Range('E27').Value = '  -4.18%  '
$ws.Range('D28').Value = '2.548.50'
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('D30').Value = '0.0₃0945'
$ws.Range('E30').Value = '  -1.00%  '
$ws.Range('D31').Value = "'" + '8.28'
$ws.Range('E31').Value = '  +0.09%  '
$ws.Range('E32').Value = '  -2.67%  '
$ws.Range('E33').Value = '  -1.61%  '
$ws.Range('E34').Value = '  -1.23%  '
$ws.Range('E35').Value = '  -2.05%  '
$ws.Range('D36').Value = "'" + '1.00'
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('D37').Value = "'" + '4.82'
$ws.Range('E37').Value = '  +0.78%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').Value = "'" + '5.58'
$ws.Range('E38').Value = '  -2.12%  '
$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D39').Value = "'" + '0.382'
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('D40').Value = "'" + '18.71'
$ws.Range('E40').Value = '  +0.18%  '
$ws.Range('D41').Value = "'" + '150.13'
$ws.Range('E41').Value = '  -1.41%  '
$ws.Range('D42').Value = "'" + '1.80'
$ws.Range('E42').Value = '  -1.49%  '
$ws.Range('D43').Value = "'" + '1.00'
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('D44').Value = "'" + '2.31'
$ws.Range('E44').Value = '  -1.07%  '
$ws.Range('D45').Value = "'" + '148.09'
$ws.Range('E45').Value = '  -1.16%  '
$ws.Range('D46').Value = "'" + '3.67'
$ws.Range('E46').Value = '  +0.26%  '
$ws.Range('D47').Value = "'" + '0.0534'
$ws.Range('E47').Value = '  -0.72%  '
$ws.Range('D48').Value = "'" + '20.25'
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('E49').Value = '  +0.12%  '
$ws.Range('D50').Value = "'" + '0.0925'
$ws.Range('E50').Value = '  +0.26%  '
$ws.Range('E51').Value = '  +0.49%  '
